$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue "D2" "27.497.03"
Set-TextValue "E2" "  -0.61%  "

Set-TextValue "D3" "1.621.41"
Set-TextValue "E3" "  -1.10%  "

Set-TextValue "E4" "  -0.01%  "

Set-TextValue "D5" "211.48"
Set-TextValue "E5" "  -0.62%  "

Set-TextValue "D6" "0.521"
Set-TextValue "E6" "  -0.55%  "

Set-TextValue "E7" "  -0.03%  "

Set-TextValue "D8" "23.16"
Set-TextValue "E8" "  +0.24%  "

Set-TextValue "D9" "0.262"
Set-TextValue "E9" "  +1.73%  "

Set-TextValue "E10" "  +0.01%  "

Set-TextValue "D11" "0.0882"
Set-TextValue "E11" "  -1.39%  "

Set-TextValue "D12" "1.851.26"

Set-TextValue "D13" "1.621.21"
Set-TextValue "E13" "  -1.16%  "

Set-TextValue "E14" "  -0.16%  "

Set-TextValue "E15" "  -1.99%  "

Set-TextValue "E16" "  +1.10%  "

Set-TextValue "D17" "27.485.51"
Set-TextValue "E17" "  -0.59%  "

Set-TextValue "D18" "229.47"
Set-TextValue "E18" "  -0.36%  "

Set-TextValue "D19" "0.0₃0717"
Set-TextValue "E19" "  -0.76%  "

Set-TextValue "D20" "7.53"
Set-TextValue "E20" "  -2.31%  "

Set-TextValue "E21" "  -0.02%  "

Set-TextValue "D22" "10.44"
Set-TextValue "E22" "  +4.09%  "

Set-TextValue "E23" "  +1.20%  "

Set-TextValue "E24" "  +8.38%  "

Set-TextValue "D25" "149.29"
Set-TextValue "E25" "  -0.28%  "

Set-TextValue "E26" "  -0.87%  "

Set-TextValue "E27" "  -0.37%  "

Set-TextValue "E28" "  +0.00%  "

Set-TextValue "E29" "  -0.78%  "

Set-TextValue "D30" "1.18"
Set-TextValue "E30" "  -0.56%  "

Set-TextValue "D31" "0.0483"
Set-TextValue "E31" "  -0.61%  "

Set-TextValue "E32" "  -0.94%  "

Set-TextValue "D33" "1.465.46"
Set-TextValue "E33" "  +1.44%  "

Set-TextValue "E34" "  -2.15%  "

Set-TextValue "E35" "  -1.58%  "

Set-TextValue "E36" "  -1.65%  "

Set-TextValue "D37" "0.944"
Set-TextValue "E37" "  +4.41%  "

Set-TextValue "E38" "  +0.27%  "

Set-TextValue "E39" "  -0.31%  "

Set-TextValue "D40" "0.553"
Set-TextValue "E40" "  -2.47%  "

Set-TextValue "E41" "  -0.02%  "

Set-TextValue "E42" "  -1.12%  "

Set-TextValue "D43" "67.84"
Set-TextValue "E43" "  -3.30%  "

Set-TextValue "E44" "  +0.46%  "

Set-TextValue "D45" "2.19"
Set-TextValue "E45" "  -2.07%  "

Set-TextValue "E46" "  -5.04%  "

Set-TextValue "E47" "  +2.79%  "

Set-TextValue "D48" "1.761.05"
Set-TextValue "E48" "  -1.19%  "

Set-TextValue "D49" "87.15"
Set-TextValue "E49" "  +1.36%  "

Set-TextValue "D50" "0.0₆0105"
Set-TextValue "E50" "  -0.54%  "

Set-TextValue "D51" "0.0993"
Set-TextValue "E51" "  +0.34%  "
